$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordering of the worker/period table in B16:G21.
# Row 16 now holds the LUIS SEGUNDO PASTRANA HERNANDEZ / period 1808 record
# (previously on row 19). The remaining JUAN DAVID MADRID OCHOA records are
# reordered with periods descending 1810, 1809, 1808, 1807, 1806, and keep
# their Valor Mora / Salario Basico values aligned to the period.

$rows = @(
    @{ Row = 16; Doc = "CC"; Id = "70526895"; Name = "LUIS SEGUNDO PASTRANA HERNANDEZ"; Period = "1808"; Mora = 42000;  Salario = 1050000 },
    @{ Row = 17; Doc = "CC"; Id = "1068391397"; Name = "JUAN DAVID MADRID OCHOA";       Period = "1810"; Mora = 31249;  Salario = 781242 },
    @{ Row = 18; Doc = "CC"; Id = "1068391397"; Name = "JUAN DAVID MADRID OCHOA";       Period = "1809"; Mora = 31249;  Salario = 781242 },
    @{ Row = 19; Doc = "CC"; Id = "1068391397"; Name = "JUAN DAVID MADRID OCHOA";       Period = "1808"; Mora = 31249;  Salario = 781242 },
    @{ Row = 20; Doc = "CC"; Id = "1068391397"; Name = "JUAN DAVID MADRID OCHOA";       Period = "1807"; Mora = 31249;  Salario = 781242 },
    @{ Row = 21; Doc = "CC"; Id = "1068391397"; Name = "JUAN DAVID MADRID OCHOA";       Period = "1806"; Mora = 27083;  Salario = 781242 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.Doc
    $ws.Cells.Item($r.Row, 3).Value = $r.Id
    $ws.Cells.Item($r.Row, 4).Value = $r.Name
    $ws.Cells.Item($r.Row, 5).Value = $r.Period
    $ws.Cells.Item($r.Row, 6).Value = $r.Mora
    $ws.Cells.Item($r.Row, 7).Value = $r.Salario
}
